$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(16, 1).Value = "2025-04-28 20:47:48"
$ws.Cells.Item(16, 2).Value = 0
